# March 24 update 3
# Adds three new columns (M: renewd, N: PlanID, O: iteration) to the
# building data sheet, with header labels in row 1 and a constant value
# per data row (rows 2-24): "after", 20180335, 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells - same bold/bordered/centered style as the other
# header cells (copy the formatting from L1, the last existing header)
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)

# New data cells for every existing data row (2 through 24)
for ($r = 2; $r -le 24; $r++) {
    $ws.Cells.Item($r, 13).Value = "after"
    $ws.Cells.Item($r, 14).Value = 20180335
    $ws.Cells.Item($r, 15).Value = 8
}
